$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in G1 and H1, matching the formatting of the existing header row
$ws.Range("G1").Value = "sauna_access"
$ws.Range("H1").Value = "steam_room"
$ws.Range("G1:H1").Font.Bold = $true
$ws.Range("G1:H1").HorizontalAlignment = -4108
$ws.Range("G1:H1").VerticalAlignment = -4160
$ws.Range("G1:H1").Borders.LineStyle = 1

# Row 2 updates
$ws.Range("C2").Value = "Group Training"
$ws.Range("E2").Value = "All Day Access"
$ws.Range("G2").Value = "7 days/week"
$ws.Range("H2").Value = "7 days/week"

# Row 3 updates
$ws.Range("C3").Value = "No Trainer"
$ws.Range("E3").Value = "All Day Access"
$ws.Range("G3").Value = "7 days/week"
$ws.Range("H3").Value = "7 days/week"

# New row 4
$ws.Range("A4").Value = "economy"
$ws.Range("B4").Value = 22222
$ws.Range("C4").Value = "Personal Trainer"
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "Morning (6 AM - 11 AM)"
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = "1 day/week"
$ws.Range("H4").Value = "1 day/week"
